# Add a new reference entry ("Adam" paper) as a new row to the
# "references" table on the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The data lives in an Excel Table ("ListObject"). Adding a ListRow grows
# the table (and its AutoFilter) by one row, same as typing data into the
# row directly beneath the table in the real workbook.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Columns, in table order: paper | publish_year | author | type | method | date_read | comment
$ws.Range("A44").Value = "ADAM: A METHOD FOR STOCHASTIC OPTIMIZATION"
$ws.Range("B44").Value = 2015
$ws.Range("C44").Value = "Diederik P. Kingma, Jimmy Lei Ba"
$ws.Range("D44").Value = "optimizer"
$ws.Range("G44").Value = "Adam optimizer"

# Leave selection on the newly-added cell, matching the saved workbook state.
$ws.Range("G44").Select()
